# week3_tues.pptx -> "added stuff for week3 thurs"
#
# 1. Slide 11 ("Equality"): merge the "see code " + "example (" runs into a
#    single "see code example (" run (leaving the SetTest.java / ")" runs,
#    and their formatting, untouched).
# 2. Slide 5 ("Lists"): merge the long "... shifting all of the elements"
#    run with the trailing "." run into one run ending in "elements."
# 3. Append a new slide 12 ("Recursion") after the last slide, using the
#    same "Title and Content" layout as the rest of the deck.

$p = $ppt.ActivePresentation

# --- 1. Slide 11: "see code " + "example (" -> "see code example (" -----
$s11 = $p.Slides.Item(11)
$body11 = $s11.Shapes.Item(2).TextFrame.TextRange
$prefixOld = "see code " + "example ("
$prefixNew = "see code example ("
$body11.Characters(1, $prefixOld.Length).Text = $prefixNew

# --- 2. Slide 5: merge trailing two runs of the LinkedList paragraph ------
$s5 = $p.Slides.Item(5)
$body5 = $s5.Shapes.Item(2).TextFrame.TextRange
$paras5 = $body5.Paragraphs()
for ($i = 1; $i -le $paras5.Count; $i++) {
    $para = $paras5.Item($i)
    if ($para.Text.StartsWith("A LinkedList provides links")) {
        $relStart = "A LinkedList".Length + 1
        $absStart = $para.Start + $relStart - 1
        $absLen = $para.Text.Length - $relStart + 1
        $merged = $para.Text.Substring($relStart - 1)
        $body5.Characters($absStart, $absLen).Text = $merged
        break
    }
}

# --- 3. New slide 12: "Recursion" -----------------------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Recursion"

$contentBody = $newSlide.Shapes.Item(2).TextFrame.TextRange
$full = "see code examples under code dir: " + "TreeRecursion.java" + " and " + "TreeNode.java"
$contentBody.Text = $full

$run1 = "see code examples under code dir: "
$run2 = "TreeRecursion.java"
$run3 = " and "
$run4 = "TreeNode.java"

$contentBody.Characters($run1.Length + 1, $run2.Length).Text = $run2
$contentBody.Characters($run1.Length + $run2.Length + $run3.Length + 1, $run4.Length).Text = $run4
